$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.704.35'
$ws.Cells.Item(2, 5).Value = '  -4.33%  '
$ws.Cells.Item(3, 4).Value = '2.981.02'
$ws.Cells.Item(3, 5).Value = '  -5.17%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '541.26'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -5.58%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '151.88'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -7.76%  '
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 5).Value = '  -1.50%  '
$ws.Cells.Item(9, 4).Value = '2.992.15'
$ws.Cells.Item(9, 5).Value = '  -5.21%  '
$ws.Cells.Item(10, 5).Value = '  -3.91%  '
$ws.Cells.Item(11, 5).Value = '  -7.23%  '
$ws.Cells.Item(12, 5).Value = '  -4.15%  '
$ws.Cells.Item(13, 4).Value = '3.504.32'
$ws.Cells.Item(13, 5).Value = '  -5.10%  '
$ws.Cells.Item(14, 5).Value = '  -2.22%  '
$ws.Cells.Item(15, 4).Value = '61.759.78'
$ws.Cells.Item(15, 5).Value = '  -4.22%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '23.92'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -4.51%  '
$ws.Cells.Item(17, 4).Value = '2.981.97'
$ws.Cells.Item(17, 5).Value = '  -5.30%  '
$ws.Cells.Item(18, 5).Value = '  -5.73%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '5.17'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.88%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.06'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -4.05%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '382.08'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -6.28%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.71'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -5.45%  '
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$ws.Cells.Item(24, 5).Value = '  -3.72%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '65.93'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -4.49%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.472'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -2.61%  '
$ws.Cells.Item(27, 4).Value = '3.103.06'
$ws.Cells.Item(27, 5).Value = '  -5.23%  '
$ws.Cells.Item(28, 5).Value = '  -2.19%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.996'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.05%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0940'
$ws.Cells.Item(30, 5).Value = '  -8.31%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '8.19'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -8.23%  '
$ws.Cells.Item(32, 5).Value = '  +0.03%  '
$ws.Cells.Item(33, 5).Value = '  -4.53%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '20.52'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.56%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '160.40'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.23%  '
$ws.Cells.Item(36, 5).Value = '  -6.13%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.92'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -6.08%  '
$ws.Cells.Item(38, 5).Value = '  -5.15%  '
$ws.Cells.Item(39, 5).Value = '  -6.81%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.55'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -8.51%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '37.58'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.87%  '
$ws.Cells.Item(42, 4).Value = '2.423.46'
$ws.Cells.Item(42, 5).Value = '  -8.33%  '
$ws.Cells.Item(43, 5).Value = '  -4.80%  '
$ws.Cells.Item(44, 5).Value = '  -7.34%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.672'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -2.94%  '
$ws.Cells.Item(46, 5).Value = '  -3.80%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '5.17'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.39%  '
$ws.Cells.Item(48, 5).Value = '  +0.13%  '
$ws.Cells.Item(49, 5).Value = '  -3.96%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0954'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.28%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '19.80'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -7.07%  '
